$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")
$ws.Activate()

$ws.Range("A2").Value = "Just a test"
$ws.Range("B2").Value = "Justtest-date"

$ws.Range("B3").Select()
